# Denmark 1st-division 2023-2024 — corrections to several match rows
# (home/away + odds + url had been mismatched against the wrong fixture)
# plus one new match row (row 96) appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowFV {
    param($row, $vals)
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 6 + $i).Value = $vals[$i]
    }
}

# --- Rows whose F:V (home..url) data was swapped/rotated between fixtures ---

Set-RowFV 6 @('Kolding IF', 2, 'Horsens', 1, 3.66, '07/07/2023 11:12', 2.86, '23/07/2023 13:32', 3.53, '07/07/2023 11:12', 3.55, '23/07/2023 13:16', 2, '07/07/2023 11:12', 2.42, '23/07/2023 13:32', 'https://www.betexplorer.com/football/denmark/1st-division/kolding-if-horsens/j1ecd7F0/')
Set-RowFV 7 @('Fredericia', 0, 'Hobro', 2, 1.81, '07/07/2023 11:12', 1.88, '23/07/2023 13:56', 3.77, '07/07/2023 11:12', 3.73, '23/07/2023 13:56', 4.15, '07/07/2023 11:12', 4.02, '23/07/2023 13:56', 'https://www.betexplorer.com/football/denmark/1st-division/fredericia-hobro/fZ1lbTpl/')

Set-RowFV 8 @('Sonderjyske', 1, 'Fredericia', 3, 1.82, '23/07/2023 14:12', 1.7, '28/07/2023 18:50', 3.92, '23/07/2023 14:12', 4.14, '28/07/2023 18:56', 3.97, '23/07/2023 14:12', 4.54, '28/07/2023 18:56', 'https://www.betexplorer.com/football/denmark/1st-division/sonderjyske-fredericia/6X7UfFeg/')
Set-RowFV 9 @('Hobro', 3, 'Naestved', 2, 2.32, '24/07/2023 03:42', 2.73, '28/07/2023 18:56', 3.35, '24/07/2023 03:42', 3.39, '28/07/2023 18:51', 3.08, '24/07/2023 03:42', 2.61, '28/07/2023 18:56', 'https://www.betexplorer.com/football/denmark/1st-division/hobro-naestved-if/ldDPeetm/')

Set-RowFV 26 @('Kolding IF', 1, 'Sonderjyske', 3, 2.08, '14/08/2023 04:12', 2.99, '18/08/2023 18:30', 3.6, '14/08/2023 04:12', 3.45, '18/08/2023 18:31', 3.38, '14/08/2023 04:12', 2.38, '18/08/2023 18:30', 'https://www.betexplorer.com/football/denmark/1st-division/kolding-if-sonderjyske/UiWf7qoG/')
Set-RowFV 28 @('Horsens', 3, 'Helsingor', 1, 2.01, '13/08/2023 22:12', 2.31, '18/08/2023 18:51', 3.79, '13/08/2023 22:12', 3.76, '18/08/2023 18:51', 3.2, '13/08/2023 22:12', 2.89, '18/08/2023 18:51', 'https://www.betexplorer.com/football/denmark/1st-division/horsens-helsingor/nTtj8PV9/')

Set-RowFV 50 @('Hobro', 0, 'Horsens', 1, 2.52, '08/09/2023 05:12', 2.42, '15/09/2023 18:52', 3.3, '08/09/2023 05:12', 3.3, '15/09/2023 18:52', 2.7, '08/09/2023 05:12', 3.05, '15/09/2023 18:52', 'https://www.betexplorer.com/football/denmark/1st-division/hobro-horsens/QLllNTcj/')
Set-RowFV 51 @('Hillerod', 0, 'B.93', 1, 2.16, '08/09/2023 05:12', 2, '15/09/2023 18:56', 3.57, '08/09/2023 05:12', 3.71, '15/09/2023 18:56', 3.21, '08/09/2023 05:12', 3.61, '15/09/2023 18:55', 'https://www.betexplorer.com/football/denmark/1st-division/hillerod-boldklubben-1893/WIuvpi53/')
Set-RowFV 52 @('Sonderjyske', 1, 'Aalborg', 1, 2.59, '08/09/2023 05:12', 2.6, '15/09/2023 18:59', 3.45, '08/09/2023 05:12', 3.57, '15/09/2023 18:59', 2.53, '08/09/2023 05:12', 2.63, '15/09/2023 18:51', 'https://www.betexplorer.com/football/denmark/1st-division/sonderjyske-aalborg/vPhhMmDd/')

Set-RowFV 75 @('Hobro', 1, 'Helsingor', 0, 1.92, '13/10/2023 18:13', 1.85, '20/10/2023 18:38', 3.68, '13/10/2023 18:13', 3.59, '20/10/2023 18:39', 3.56, '13/10/2023 18:13', 4.38, '20/10/2023 18:36', 'https://www.betexplorer.com/football/denmark/1st-division/hobro-helsingor/YqNTnExr/')
Set-RowFV 76 @('Hillerod', 1, 'Fredericia', 1, 2.84, '13/10/2023 18:13', 3.03, '20/10/2023 18:58', 3.6, '13/10/2023 18:13', 3.79, '20/10/2023 18:58', 2.26, '13/10/2023 18:13', 2.31, '20/10/2023 18:35', 'https://www.betexplorer.com/football/denmark/1st-division/hillerod-fredericia/G4z2ifUQ/')

Set-RowFV 79 @('B.93', 0, 'Sonderjyske', 4, 5.01, '22/10/2023 15:12', 8.03, '27/10/2023 18:58', 4.47, '22/10/2023 15:12', 5.46, '27/10/2023 18:58', 1.57, '22/10/2023 15:12', 1.34, '27/10/2023 18:50', 'https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-sonderjyske/2PEpqWy7/')
Set-RowFV 80 @('Hobro', 2, 'Koge', 1, 1.74, '22/10/2023 16:12', 1.65, '27/10/2023 18:51', 3.92, '22/10/2023 16:12', 4.26, '27/10/2023 18:51', 4.42, '22/10/2023 16:12', 4.85, '27/10/2023 18:51', 'https://www.betexplorer.com/football/denmark/1st-division/hobro-koge/UBGxoh7f/')
Set-RowFV 81 @('Horsens', 1, 'Hillerod', 3, 1.84, '20/10/2023 19:12', 2.17, '27/10/2023 18:51', 3.9, '20/10/2023 19:12', 3.59, '27/10/2023 18:51', 3.63, '20/10/2023 19:12', 3.27, '27/10/2023 18:51', 'https://www.betexplorer.com/football/denmark/1st-division/horsens-hillerod/G4dX6ZTD/')

# --- New row 96: append one more fixture, copying row 95's formatting first ---

$ws.Range("A95:V95").Copy()
$ws.Range("A96:V96").PasteSpecial(-4122)

$ws.Cells.Item(96,1).Value = 95
$ws.Cells.Item(96,2).Value = "denmark"
$ws.Cells.Item(96,3).Value = "1st-division"
$ws.Cells.Item(96,4).Value = "2023-2024"
$ws.Cells.Item(96,5).Value = 45243.79166666666
$ws.Cells.Item(96,5).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Set-RowFV 96 @('Hobro', 0, 'Hillerod', 0, 2.31, '06/11/2023 19:12', 2.32, '13/11/2023 18:58', 3.62, '06/11/2023 19:12', 3.33, '13/11/2023 18:58', 2.76, '06/11/2023 19:12', 3.2, '13/11/2023 18:58', 'https://www.betexplorer.com/football/denmark/1st-division/hobro-hillerod/MPau14ch/')
